$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements (row => new text), single-column table.
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "201"
    5  = "0.00003"
    6  = "0.00053"
    7  = "0.00017"
    8  = "0.00005"
    9  = "0.00025"
    10 = "0.00035"
    11 = "0.00040"
    12 = "0.03942"
    44 = "99.95"
    45 = "0.04"
    46 = "77"
}

foreach ($row in $updates.Keys) {
    $t.Cell($row, 1).Range.Text = $updates[$row]
}
